# docs/stmts.xlsx update: "updated with work from last night and this morning"
#
# TODO sheet (sheet1): drop the "Example" column (D) entirely, rename the
# Status header, mark "For loops" Done, re-point the active view.
# Keywords sheet (sheet2): fill in the new "NEXT" row (BTOKEN_NEXT), and
# swap which sheet/cell is active/selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # TODO
$ws2 = $wb.Worksheets.Item(2)   # Keywords

# --- TODO sheet -----------------------------------------------------------

# Header: "Statis" -> "Status"; drop the "Example" header (D1)
$ws1.Range("C1").Value = "Status"
$ws1.Range("D1").ClearContents()

# Row 2 (User defined functions): drop the worked example (D2)
$ws1.Range("D2").ClearContents()

# Row 3 (For loops) is now Done too
$ws1.Range("C3").Value = "Done"

# Row 7 (IF/THEN): drop the worked example (D7)
$ws1.Range("D7").ClearContents()

# --- Keywords sheet ---------------------------------------------------------

# Row 15: NEXT keyword row gets filled in like the other token rows
$ws2.Range("C15").Value = "X"
$ws2.Range("D15").Value = "X"
$ws2.Range("E15").Value = "X"
$ws2.Range("F15").Value = "X"
$ws2.Range("G15").Value = "BTOKEN_NEXT"

# --- View state -------------------------------------------------------------
# Active sheet flips from Keywords back to TODO, with new selections on each.
$null = $ws2.Range("H15").Select()
$ws1.Activate()
$null = $ws1.Range("C2").Select()
